# This script reproduces the authored edit to streamlit_apps.xlsx:
#   1. Cell B2's text gains an "https://" scheme prefix
#      ("santi-nue.github.io/" -> "https://santi-nue.github.io/").
#   2. A hyperlink pointing at that same URL is attached to cell B2
#      (it did not have one before).
#   3. The worksheet's active cell / selection moves from A3 to B2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newUrl = "https://santi-nue.github.io/"

# 1. Update B2's displayed text to include the "https://" prefix.
$ws.Range("B2").Value = $newUrl

# 2. Attach a new hyperlink to B2 that targets the same URL. The
#    TextToDisplay argument mirrors what the source workbook stores
#    for this hyperlink's cached display text.
$ws.Hyperlinks.Add($ws.Range("B2"), $newUrl, "", "", "https://santi") | Out-Null

# Adding the hyperlink can resync the cell's text to the TextToDisplay
# value above, so re-assert the full URL as the cell's actual content.
$ws.Range("B2").Value = $newUrl

# 3. Select B2 so it becomes the active cell in the saved view.
$ws.Range("B2").Select() | Out-Null
